$wb = $excel.ActiveWorkbook

# --- "Anthony Davis" sheet: add a new game row (row 8) and push the
#     averages row down to row 9, expanding the AVERAGE ranges. ---
$ws1 = $wb.Worksheets.Item("Anthony Davis")
$ws1.Rows.Item(8).Insert()

$ws1.Range("A8").Value = 363
$ws1.Range("B8").Value = 43
$ws1.Range("C8").Value = 456
$ws1.Range("D8").Value = 12
$ws1.Range("E8").Value = 781
$ws1.Range("F8").Value = 79.61
$ws1.Range("G8").Value = 84.32
$ws1.Range("H8").Value = 84.66
$ws1.Range("I8").Value = 782

$ws1.Range("A9").Formula = "=AVERAGE(A2:A8)"
$ws1.Range("B9").Formula = "=AVERAGE(B2:B8)"
$ws1.Range("C9").Formula = "=AVERAGE(C2:C8)"
$ws1.Range("D9").Formula = "=AVERAGE(D2:D8)"
$ws1.Range("E9").Formula = "=AVERAGE(E2:E8)"
$ws1.Range("F9").Formula = "=AVERAGE(F2:F8)"
$ws1.Range("G9").Formula = "=AVERAGE(G2:G8)"
$ws1.Range("H9").Formula = "=AVERAGE(H2:H8)"
$ws1.Range("I9").Formula = "=AVERAGE(I2:I8)"

# --- "D Angelo Russell" sheet: add three new game rows (9, 10, 11) and
#     push the averages row down to row 12, expanding the AVERAGE ranges. ---
$ws2 = $wb.Worksheets.Item("D Angelo Russell")
$ws2.Rows.Item(9).Insert()
$ws2.Rows.Item(9).Insert()
$ws2.Rows.Item(9).Insert()

$ws2.Range("A9").Value = 150
$ws2.Range("B9").Value = 5
$ws2.Range("C9").Value = 330
$ws2.Range("D9").Value = 25
$ws2.Range("E9").Value = 330
$ws2.Range("F9").Value = 45.45
$ws2.Range("G9").Value = 46.21
$ws2.Range("H9").Value = 48.39
$ws2.Range("I9").Value = 202

$ws2.Range("A10").Value = 79
$ws2.Range("B10").Value = 45
$ws2.Range("C10").Value = 140
$ws2.Range("D10").Value = 20
$ws2.Range("E10").Value = 215
$ws2.Range("F10").Value = 56.43
$ws2.Range("G10").Value = 72.5
$ws2.Range("H10").Value = 72.24
$ws2.Range("I10").Value = 186

$ws2.Range("A11").Value = 79
$ws2.Range("B11").Value = 45
$ws2.Range("C11").Value = 140
$ws2.Range("D11").Value = 20
$ws2.Range("E11").Value = 215
$ws2.Range("F11").Value = 56.43
$ws2.Range("G11").Value = 72.5
$ws2.Range("H11").Value = 72.24
$ws2.Range("I11").Value = 186

$ws2.Range("A12").Formula = "=AVERAGE(A2:A11)"
$ws2.Range("B12").Formula = "=AVERAGE(B2:B11)"
$ws2.Range("C12").Formula = "=AVERAGE(C2:C11)"
$ws2.Range("D12").Formula = "=AVERAGE(D2:D11)"
$ws2.Range("E12").Formula = "=AVERAGE(E2:E11)"
$ws2.Range("F12").Formula = "=AVERAGE(F2:F11)"
$ws2.Range("G12").Formula = "=AVERAGE(G2:G11)"
$ws2.Range("H12").Formula = "=AVERAGE(H2:H11)"
$ws2.Range("I12").Formula = "=AVERAGE(I2:I11)"

# The "final" sheet's cross-sheet formulas (rows 3 and 5) referencing the
# averages rows on "Anthony Davis" and "D Angelo Russell" are automatically
# re-pointed by Excel's row-insert logic (A8->A9, A9->A12, etc.), so no
# further action is required there.
